$d = $word.ActiveDocument

$replacements = @(
    @("749÷6=", "975÷7="),
    @("843÷5=", "743÷9="),
    @("463÷3=", "829÷6="),
    @("373÷6=", "776÷3="),
    @("667÷2=", "520÷7="),
    @("896÷2=", "239÷6="),
    @("305÷3=", "721÷4="),
    @("169÷5=", "985÷5="),
    @("753÷6=", "715÷4="),
    @("830÷3=", "293÷7="),
    @("419÷5=", "856÷6="),
    @("544÷7=", "471÷8="),
    @("457÷2=", "367÷4="),
    @("299÷7=", "280÷9="),
    @("544÷6=", "891÷6="),
    @("119÷9=", "197÷4="),
    @("988÷9=", "232÷2="),
    @("462÷5=", "926÷3="),
    @("183÷9=", "546÷2="),
    @("961÷7=", "784÷9="),
    @("494÷4=", "487÷7="),
    @("567÷2=", "641÷9="),
    @("783÷5=", "692÷2="),
    @("102÷8=", "451÷4="),
    @("230÷7=", "133÷2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
